$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.282.61"
$ws.Range("E2").Value = "'  +1.46%  "
$ws.Range("D3").Value = "'2.655.06"
$ws.Range("E3").Value = "'  -0.17%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'609.03"
$ws.Range("E5").Value = "'  -0.51%  "
$ws.Range("D6").Value = "'148.84"
$ws.Range("E6").Value = "'  +3.60%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "'  +0.66%  "
$ws.Range("E9").Value = "'  +2.25%  "
$ws.Range("D10").Value = "'0.387"
$ws.Range("E10").Value = "'  +7.07%  "
$ws.Range("D11").Value = "'5.62"
$ws.Range("E11").Value = "'  +0.15%  "
$ws.Range("E12").Value = "'  -0.95%  "
$ws.Range("D13").Value = "'27.65"
$ws.Range("E13").Value = "'  +1.36%  "
$ws.Range("D14").Value = "'3.130.78"
$ws.Range("E14").Value = "'  -0.21%  "
$ws.Range("D15").Value = "'64.132.44"
$ws.Range("E15").Value = "'  +1.48%  "
$ws.Range("E16").Value = "'  +2.50%  "
$ws.Range("D17").Value = "'2.651.38"
$ws.Range("E17").Value = "'  -0.62%  "
$ws.Range("D18").Value = "'12.03"
$ws.Range("E18").Value = "'  +5.12%  "
$ws.Range("D19").Value = "'4.60"
$ws.Range("E19").Value = "'  +4.04%  "
$ws.Range("D20").Value = "'346.89"
$ws.Range("E20").Value = "'  +1.59%  "
$ws.Range("D21").Value = "'6.93"
$ws.Range("E21").Value = "'  +0.92%  "
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "'  -0.48%  "
$ws.Range("D24").Value = "'66.37"
$ws.Range("E24").Value = "'  -0.93%  "
$ws.Range("D25").Value = "'1.67"
$ws.Range("E25").Value = "'  +9.06%  "
$ws.Range("E26").Value = "'  +4.44%  "
$ws.Range("D27").Value = "'9.40"
$ws.Range("E27").Value = "'  +8.90%  "
$ws.Range("D28").Value = "'560.02"
$ws.Range("E28").Value = "'  +2.77%  "
$ws.Range("D29").Value = "'8.19"
$ws.Range("E29").Value = "'  +4.76%  "
$ws.Range("D30").Value = "'0.162"
$ws.Range("E30").Value = "'  -1.35%  "
$ws.Range("E31").Value = "'  -0.07%  "
$ws.Range("E32").Value = "'  +1.09%  "
$ws.Range("D33").Value = "'0.0₃0852"
$ws.Range("E33").Value = "'  +5.87%  "
$ws.Range("E34").Value = "'  -1.39%  "
$ws.Range("D35").Value = "'5.35"
$ws.Range("E35").Value = "'  +3.36%  "
$ws.Range("D36").Value = "'168.58"
$ws.Range("E36").Value = "'  -2.82%  "
$ws.Range("D37").Value = "'0.408"
$ws.Range("E37").Value = "'  +0.36%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "'  -0.05%  "
$ws.Range("E39").Value = "'  +5.25%  "
$ws.Range("D40").Value = "'19.38"
$ws.Range("E40").Value = "'  +1.21%  "
$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("D42").Value = "'167.43"
$ws.Range("E42").Value = "'  -4.85%  "
$ws.Range("D43").Value = "'40.46"
$ws.Range("E43").Value = "'  +0.84%  "
$ws.Range("D44").Value = "'3.85"
$ws.Range("E44").Value = "'  +2.50%  "
$ws.Range("D45").Value = "'22.20"
$ws.Range("E45").Value = "'  +0.09%  "
$ws.Range("D46").Value = "'0.0572"
$ws.Range("E46").Value = "'  -0.67%  "
$ws.Range("D47").Value = "'0.630"
$ws.Range("E47").Value = "'  -0.35%  "
$ws.Range("B48").Value = "'VeChain"
$ws.Range("C48").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0246"
$ws.Range("E48").Value = "'  +2.29%  "
$ws.Range("B49").Value = "'dogwifhat"
$ws.Range("C49").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'1.99"
$ws.Range("E49").Value = "'  +14.31%  "
$ws.Range("D50").Value = "'0.0963"
$ws.Range("E50").Value = "'  +0.18%  "
$ws.Range("D51").Value = "'19.01"
$ws.Range("E51").Value = "'  +1.69%  "
